$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# 1. Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 00:59"

# 2. Update Estados Unidos row (row 4) statistics
$ws.Range("B4").Value = 123311
$ws.Range("C4").Value = 19185
$ws.Range("E4").Value = 117869
$ws.Range("G4").Value = 515
$ws.Range("H4").Value = 2211

# 3. Update Canada row (row 18) statistics
$ws.Range("B18").Value = 5655
$ws.Range("C18").Value = 898
$ws.Range("E18").Value = 5199

# 4. Reorder Uruguay / Bosnia y Herzegovina and refresh their data.
#    Row 79 becomes Bosnia y Herzegovina (with its updated stats),
#    row 80 becomes Uruguay (keeping its previous, unchanged stats).
$ws.Range("A79").Value = "Bosnia y Herzegovina"
$ws.Range("B79").Value = 278
$ws.Range("C79").Value = 41
$ws.Range("D79").Value = 8
$ws.Range("E79").Value = 264
$ws.Range("F79").Value = 1
$ws.Range("G79").Value = 2
$ws.Range("H79").Value = 6

$ws.Range("A80").Value = "Uruguay"
$ws.Range("B80").Value = 274
$ws.Range("C80").Value = 36
$ws.Range("D80").Value = 0
$ws.Range("E80").Value = 274
$ws.Range("F80").Value = 8
$ws.Range("G80").Value = 0
$ws.Range("H80").Value = 0
